# Update "想去人数" (want-to-go count) figures on the 展览, 演出, and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 2035
$wsExhibit.Range("F5").Value = 345
$wsExhibit.Range("F8").Value = 13320
$wsExhibit.Range("F10").Value = 48
$wsExhibit.Range("F11").Value = 5394
$wsExhibit.Range("F12").Value = 561
$wsExhibit.Range("F15").Value = 44
$wsExhibit.Range("F18").Value = 147
$wsExhibit.Range("F21").Value = 7383
$wsExhibit.Range("F23").Value = 3665

$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 31

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 2035
$wsAll.Range("F5").Value = 345
$wsAll.Range("F6").Value = 31
$wsAll.Range("F9").Value = 13320
$wsAll.Range("F11").Value = 48
$wsAll.Range("F12").Value = 5394
$wsAll.Range("F13").Value = 561
$wsAll.Range("F16").Value = 44
$wsAll.Range("F19").Value = 147
$wsAll.Range("F23").Value = 7383
$wsAll.Range("F25").Value = 3665
